$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Range("A4").Value = 10090.120000000001
$ws.Range("B4").Value = 10073
$ws.Range("C4").Value = 309.02999999999997
$ws.Range("D4").Value = 309.55
$ws.Range("E4").Value = $false
$ws.Range("F4").Value = 0.17
$ws.Range("G4").Value = 42608.624780092592
$ws.Range("G4").NumberFormat = "m/d/yy h:mm"
$ws.Range("H4").Value = $true

# Row 5
$ws.Range("A5").Value = 10107.27
$ws.Range("B5").Value = 10090.120000000001
$ws.Range("C5").Value = 309.02999999999997
$ws.Range("D5").Value = 309.55
$ws.Range("E5").Value = $false
$ws.Range("F5").Value = 0.17
$ws.Range("G5").Value = 42608.63784722222
$ws.Range("G5").NumberFormat = "m/d/yy h:mm"
$ws.Range("H5").Value = $true
